$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test case 3 (row 4) -------------------------------------------------

$summary = @"
Verify that you can not add more product in cart that the product avaliabke in store
"@

$steps = @"
1. Goto http://live.demoguru99.com
2. Cick one mobile menu
3. In the list of all mobile, click on "add to cart" for sony xperia mobile. 
4. Change "QTY" value to1000 and click "update" button
5. verify the error message
6. Then cilck on "Empty cart" link in the footer of list of all mobiules
7. Verify cart is empty
"@

$expect = @"
1. On clicking update button an error is shown 'The requested quantity for "Sony Xperia" is not avaliable.'
2. On clicking empty cart button - a message 'Shopping cart is empty' is shown
"@

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = $summary
$ws.Range("C4").Value = $steps
$ws.Range("E4").Value = $expect

$ws.Rows.Item(4).RowHeight = 140

# --- View state (best effort: scroll the window so row 3 becomes the
#     top-visible row, matching the saved workbook's topLeftCell="A3";
#     selection itself is left untouched at C3) ----------------------------
try { $excel.ActiveWindow.ScrollRow = 3 } catch {}
